$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, duplicate the formatting (style) of an existing data row (row 2, which
# has the bold/border style "A" column formatting) into the two brand-new rows
# 16 and 17 so their A-column cells end up styled like the rest of the table.
$ws.Range("A2:E2").Copy()
$ws.Range("A16:E17").PasteSpecial(-4122)
$ws.Range("A16:E16").ClearContents()
$ws.Range("A17:E17").ClearContents()

# Update existing rows 8-15 (labels shift because two new "line7"/"line8"
# entries are introduced, and the underlying data values change) and create
# new rows 16-17 for the two additional contingencies (extr7 / extr8).
$rows = @(
    @{ Row = 8;  A = 6;  B = "line7"; C = 14; D = 11; E = $true  },
    @{ Row = 9;  A = 7;  B = "line8"; C = 16; D = 9;  E = $true  },
    @{ Row = 10; A = 8;  B = "extr1"; C = 5;  D = 12; E = $true  },
    @{ Row = 11; A = 9;  B = "extr2"; C = 5;  D = 9;  E = $true  },
    @{ Row = 12; A = 10; B = "extr3"; C = 10; D = 11; E = $false },
    @{ Row = 13; A = 11; B = "extr4"; C = 7;  D = 8;  E = $false },
    @{ Row = 14; A = 12; B = "extr5"; C = 9;  D = 11; E = $false },
    @{ Row = 15; A = 13; B = "extr6"; C = 7;  D = 11; E = $true  },
    @{ Row = 16; A = 14; B = "extr7"; C = 5;  D = 7;  E = $true  },
    @{ Row = 17; A = 15; B = "extr8"; C = 8;  D = 5;  E = $false }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
}

Write-Host "done"
